$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix C454: clear the old "NA" text, leaving an empty (but still
# text-typed) cell, matching the rest of the "Rien ne nous concerne
# aujourd'hui !" rows (e.g. C453) which store an empty inline string. ---
$ws.Range("C454").Value = "'"
$ws.Range("C454").Style = "Normal"

# --- Append the new rows (455-466) produced by the latest script run. ---
$newRows = @(
    @("2026-02-11", "agriculture biologique", 4, 1),
    @("2026-02-11", "agriculture biologique", 47, 1),
    @("2026-02-11", "eaux souterraines", 49, 1),
    @("2026-02-11", "ruissellement", 54, 1),
    @("2026-02-11", "herbicides", 60, 1),
    @("2026-02-11", "eaux de surface", 60, 2),
    @("2026-02-11", "phytolicence", 73, 1),
    @("2026-02-11", "herbicides", 73, 1),
    @("2026-02-11", "eaux souterraines", 73, 1),
    @("2026-02-11", "effluents phytopharmaceutiques", 73, 1),
    @("2026-02-11", "agriculture biologique", 172, 4),
    @("2026-02-11", "bonnes pratiques", 183, 1)
)

$startRow = 455
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Column A holds the date as plain text ("2026-02-11"), not a real
    # date value, so force text with a leading apostrophe and strip the
    # resulting quote-prefix style back off afterwards.
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 1).Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
